# Timesheet update: log hours for 2024-02-12 (row 33, previously blank) and
# extend the sheet with day-rows for 2024-02-13 .. 2024-02-15 (rows 34-36,
# previously placeholder/blank rows) plus one more templated blank row
# (row 37), mirroring the existing day-row layout (date in col A, per-LO
# minute formulas in B:E, a SUM total in F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date column (A) uses a short-date custom number format; apply it to
# every row we are about to populate so the new date cells render like the
# existing ones (reuses the sheet's existing date style instead of minting
# a near-duplicate one).
$ws.Range("A33:A37").NumberFormat = "[$-409]d\-mmm;@"

# --- Row 33: 2024-02-12 now has logged time against LO1 and LO2 ---
$ws.Range("A33").Value = 45334
$ws.Range("B33").Formula = "=(1/60)*(13)"
$ws.Range("C33").Formula = "=(1/60)*(21+8)"

# --- Rows 34-36: new day rows for 2024-02-13 / 14 / 15, no time logged yet ---
$ws.Range("A34").Value = 45335
$ws.Range("A35").Value = 45336
$ws.Range("A36").Value = 45337

$ws.Range("B34:E36").Formula = "=(1/60)*(0)"
$ws.Range("F34").Formula = "=SUM(B34:E34)"
$ws.Range("F35").Formula = "=SUM(B35:E35)"
$ws.Range("F36").Formula = "=SUM(B36:E36)"

# --- Row 37: one more templated blank row, date left empty ---
$ws.Range("B37:E37").Formula = "=(1/60)*(0)"
$ws.Range("F37").Formula = "=SUM(B37:E37)"
$ws.Range("B37:F37").NumberFormat = "0.00"

# --- Selection: user had clicked into B33 after entering the data, with the
#     view scrolled down a bit so row 33 is visible near the top ---
$ws.Range("B33").Select()
